# Update "想去人数" (interested-attendee count) values in column F
# across the 展览 (Exhibition), 演出 (Performance) and 全部类型
# (All types) sheets, matching the refreshed scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1897
$ws.Range("F3").Value = 26
$ws.Range("F4").Value = 862
$ws.Range("F5").Value = 110
$ws.Range("F6").Value = 43
$ws.Range("F8").Value = 251
$ws.Range("F10").Value = 156
$ws.Range("F11").Value = 144
$ws.Range("F13").Value = 4464
$ws.Range("F15").Value = 36
$ws.Range("F16").Value = 489
$ws.Range("F17").Value = 439
$ws.Range("F18").Value = 12
$ws.Range("F20").Value = 1114
$ws.Range("F21").Value = 2207
$ws.Range("F22").Value = 374
$ws.Range("F23").Value = 60
$ws.Range("F24").Value = 40
$ws.Range("F26").Value = 2174
$ws.Range("F27").Value = 84
$ws.Range("F28").Value = 68
$ws.Range("F30").Value = 153
$ws.Range("F31").Value = 100
$ws.Range("F33").Value = 217
$ws.Range("F34").Value = 33

# --- Sheet 2: 演出 (Performance) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 36

# --- Sheet 4: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1897
$ws.Range("F3").Value = 26
$ws.Range("F4").Value = 862
$ws.Range("F5").Value = 110
$ws.Range("F6").Value = 43
$ws.Range("F8").Value = 251
$ws.Range("F10").Value = 156
$ws.Range("F11").Value = 144
$ws.Range("F13").Value = 36
$ws.Range("F14").Value = 4464
$ws.Range("F16").Value = 36
$ws.Range("F17").Value = 489
$ws.Range("F18").Value = 439
$ws.Range("F19").Value = 12
$ws.Range("F21").Value = 1114
$ws.Range("F22").Value = 2207
$ws.Range("F23").Value = 374
$ws.Range("F24").Value = 60
$ws.Range("F25").Value = 40
$ws.Range("F27").Value = 2174
$ws.Range("F28").Value = 84
$ws.Range("F29").Value = 68
$ws.Range("F31").Value = 153
$ws.Range("F32").Value = 100
$ws.Range("F34").Value = 217
$ws.Range("F35").Value = 33
